$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 24 de Abril de 2020 a las 00:22"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 878006
$ws.Range("C4").Value = 29289
$ws.Range("D4").Value = 85162
$ws.Range("E4").Value = 743107
$ws.Range("G4").Value = 2078
$ws.Range("H4").Value = 49737

# Row 16 - Canada
$ws.Range("B16").Value = 42081
$ws.Range("C16").Value = 1891
$ws.Range("D16").Value = 14748
$ws.Range("E16").Value = 25190
$ws.Range("G16").Value = 169
$ws.Range("H16").Value = 2143

# Row 50 - Colombia
$ws.Range("B50").Value = 4561
$ws.Range("C50").Value = 205
$ws.Range("D50").Value = 927
$ws.Range("E50").Value = 3419
$ws.Range("G50").Value = 9
$ws.Range("H50").Value = 215

# Row 88 - Tunez
$ws.Range("B88").Value = 918
$ws.Range("C88").Value = 9
$ws.Range("E88").Value = 690
